$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group header / summary rows (2 and 3)
$ws.Range("A2").Value = "Zipper (1)"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "274000.0"
$ws.Range("A3").Value = "    BDT (1)"
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "274000.0"

# Row 4 updates (the Purchase Order record)
$ws.Range("A4").Value = "Zipper"
$ws.Range("B4").Value = "Shamsuddin Ahamed"
$ws.Range("C4").Value = 45806.60496527778
$ws.Range("D4").Value = "BDT"
$ws.Range("F4").Value = "MD. MONIR HOSSAIN"
$ws.Range("G4").Value = "Md. Shahid Hossain"
$ws.Range("H4").Value = 45806.72810185186
$ws.Range("I4").Value = "P18037"
$ws.Range("L4").Value = "Normal"
$ws.Range("M4").Value = "Dyeing CIP (Mr Anup)"
$ws.Range("N4").Value = "RFQ Sent"
$ws.Range("O4").Value = 274000
$ws.Range("P4").Value = "Vision Tex CO"
$ws.Range("Q4").Value = "BDT"
$ws.Range("R4").Value = ""
$ws.Range("S4").Value = ""
$ws.Range("T4").Value = "By Road"
